# Classes for working with xls files were added.
#
# The only semantically meaningful change to this workbook fixture is that
# the second worksheet ("Лист2") is renamed to "My parent's accounting".
#
# (The rest of the original diff - the x15ac:absPath session path, the
# xr:revisionPtr documentId GUID, the bookViews window geometry, and the
# internal re-ordering of the <mergeCells> entries on that sheet - are all
# byproducts of Excel re-saving the file on the author's machine: a local
# file-system path, freshly generated revision/session GUIDs, the OS
# window position/size at save time, and the iteration order Excel used
# internally when rewriting the merged-cell table. None of that reflects
# an actual edit to the workbook's data/structure, and these values are
# not meant to be deterministically reproduced here - the mergeCells
# *set* itself is unchanged, only its incidental on-disk ordering differs.)

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

$ws2.Name = "My parent's accounting"
